$p = $ppt.ActivePresentation
$s = $p.Slides.Item(31)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Locate the "Pacman:" paragraph (3rd paragraph in this text box).
$pacmanPara = $tr.Paragraphs(3, 1)

# Update the label run to include a trailing space, matching the
# "Spaceinvader: " / "Tic Tac Toe: " runs above it.
$labelRun = $tr.Characters($pacmanPara.Start, 7)
$labelRun.Text = "Pacman: "

$pacmanPara = $tr.Paragraphs(3, 1)
$urlText = "https://github.com/daleharvey/pacman"
$urlLen = $urlText.Length

# Insert the URL text, plus a trailing sentinel character. Using the
# sentinel keeps the new hyperlink run from being the very last run in
# the paragraph while we set its ActionSettings hyperlink, which avoids
# the host re-syncing a stray <a:endParaRPr> onto the paragraph. We
# delete the sentinel afterwards.
$pacmanPara.InsertAfter($urlText + "Z") | Out-Null

$urlStart = $pacmanPara.Start + 8
$urlRun = $tr.Characters($urlStart, $urlLen)

$urlRun.Font.Underline = $true
$urlRun.Font.Color.ObjectThemeColor = 11  # msoThemeColorHyperlink

$action = $urlRun.ActionSettings(1)  # ppMouseClick
$action.Hyperlink.Address = $urlText

$sentinel = $tr.Characters($urlStart + $urlLen, 1)
$sentinel.Text = ""
